$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 138; this shifts the existing rows 138-170
# down to 139-171 (preserving all of their data/formatting).
$ws.Rows(138).Insert()

# Populate the newly inserted row 138 with the new data record.
$ws.Cells.Item(138, 1).Value = 7
$ws.Cells.Item(138, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(138, 3).Value = "Ñuble"
$ws.Cells.Item(138, 4).Value = 45211
$ws.Cells.Item(138, 5).Value = 16
$ws.Cells.Item(138, 6).Value = 100112031
$ws.Cells.Item(138, 7).Value = "Poroto verde"
$ws.Cells.Item(138, 8).Value = "Sin especificar"
$ws.Cells.Item(138, 9).Value = "Primera"
$ws.Cells.Item(138, 10).Value = 30
$ws.Cells.Item(138, 11).Value = 29000
$ws.Cells.Item(138, 12).Value = 29000
$ws.Cells.Item(138, 13).Value = 29000
$ws.Cells.Item(138, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(138, 15).Value = "Perú"
$ws.Cells.Item(138, 16).Value = 1160
$ws.Cells.Item(138, 17).Value = 25
$ws.Cells.Item(138, 18).Value = "Hortaliza"
